$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Remove the now-unused trailing columns (X:AG) that existed on rows 1-19.
#    After the edit the sheet only spans columns A:W.
# ---------------------------------------------------------------------------
$ws.Range("X1:AG19").Delete()

# ---------------------------------------------------------------------------
# 2. Update the row 2 header labels (C2:W2) to reflect the new HKL ordering.
# ---------------------------------------------------------------------------
$headerLabels = "[4, 2, 2]|[5, 1, 1]|[2, 2, 2]|[1, 1, 1]|[3, 1, 1]|[3, 3, 1]|[2, 2, 0]|[2, 0, 0]|[3, 3, 3]|[4, 0, 0]|[4, 2, 0]|1Pair-A|1Pair-B|2Pairs-A|2Pairs-B|3Pairs-A|3Pairs-B|3Pairs-C|4Pairs|5A4F|MaxUnique" -split "\|"
$col = 3
foreach ($lbl in $headerLabels) {
    $ws.Cells.Item(2, $col).Value = $lbl
    $col++
}

# ---------------------------------------------------------------------------
# 3. Rewrite the data rows (3-23): index column A, scheme-name column B, and
#    the 21 numeric average-HW ratios in columns C:W. Rows 16-19 now hold the
#    new "Holden" simulation results, and the HexGrid rows have moved down
#    to rows 20-23.
# ---------------------------------------------------------------------------
$dataRows = @"
3|1|BT8Hex_2.5|0.9992379967779946,1.001974158938575,0.997966211717362,0.997966211717362,1.000652425766276,0.9988742210910363,0.9992473853321913,1.00303459963446,0.997966211717362,1.00303459963446,1.000610776919891,0.997966211717362,0.9992473853321913,1.001140992483326,0.9999499055492336,1.000082732228004,1.000978136910976,1.000082732228005,1.000225155612573,0.9997733668335304,1.000199722022223
4|2|BT8Hex_5|0.9985289825605073,1.003811103022706,0.9960738592919095,0.9960738592919095,1.001259503391195,0.9978266757183147,0.9985470258190545,1.00585829313776,0.9960738592919095,1.00585829313776,1.001179075508591,0.9960738592919095,0.9985470258190545,1.002202659478407,0.9999032646051249,1.000159726082908,1.001888274116003,1.000159726082908,1.00043467040998,0.9995625081863657,1.000385564806255
5|3|BT8Hex_10|0.9971813481762857,1.007303800168198,0.9924773061108922,0.9924773061108922,1.002413774801749,0.9958347184110204,0.9972143877006703,1.011227412354476,0.9924773061108922,1.011227412354476,1.002259070501342,0.9924773061108922,0.9972143877006703,1.004220900027573,0.9998140812512094,1.000306368722013,1.003618524952298,1.000306368722013,1.000833220241947,0.9991620374157358,1.000738977278079
6|4|BT8Hex_15|0.9958661300284649,1.010712512833168,0.9889672931603133,0.9889672931603133,1.003540291746927,0.9938906578981959,0.9959137608553623,1.016467418837953,0.9889672931603133,1.016467418837953,1.003313073073102,0.9889672931603133,0.9959137608553623,1.006190589846658,0.9997270263011446,1.000449490951209,1.005307157146747,1.000449490951209,1.001222191150139,0.9987712115521739,1.001083892304186
7|5|Spiral2.5|0.999941679871321,1.000152423442629,0.9998446626474541,0.9998446626474541,1.000050370232257,0.9999128687601375,0.999940769174183,1.00023457939932,0.9998446626474541,1.00023457939932,1.000046535700021,0.9998446626474541,0.999940769174183,1.000087674286751,0.9999955697032201,1.000006670406986,1.00007523960192,1.000006670406986,1.000017595363304,0.9999830088201337,1.000015486153415
8|6|Spiral5|0.9998403779144506,1.00041710852302,0.9995748074702432,0.9995748074702432,1.000137838827561,0.9997615696652529,0.9998379790575023,1.000641913111906,0.9995748074702432,1.000641913111906,1.00012739038042,0.9995748074702432,0.9998379790575023,1.000239946084704,0.9999879089425314,1.000018233213217,1.000205910332323,1.000018233213217,1.000048134616803,0.9999534691874912,1.000042373118795
9|7|Spiral7.5|0.9997686316078427,1.000604574171294,0.9993836868234174,0.9993836868234174,1.000199789651228,0.9996544096891047,0.9997651724452485,1.000930411754366,0.9993836868234174,1.000930411754366,1.00018465373781,0.9993836868234174,0.9997651724452485,1.000347792099808,0.9999824810482384,1.000026423674344,1.000298457950281,1.000026423674344,1.000069765168565,0.9999325494995356,1.000061416235039
10|8|Spiral10|0.9995087933871268,1.001283454917904,0.998691503835595,0.998691503835595,1.000424134506672,0.9992663530223266,0.9995015561672,1.001975161694816,0.998691503835595,1.001975161694816,1.000392048898538,0.998691503835595,0.9995015561672,1.000738358931008,0.999962845336936,1.000056073899204,1.000633617456229,1.000056073899204,1.000148089051071,0.9998567720079755,1.000130375803772
11|9|Spiral15|0.9991626369816176,1.00218821658437,0.9977694649724077,0.9977694649724077,1.000723125728503,0.9987491270983949,0.9991499296551178,1.003367596988901,0.9977694649724077,1.003367596988901,1.00066828750752,0.9977694649724077,0.9991499296551178,1.001258763322009,0.9999365276918102,1.000095663872142,1.001080217457507,1.000095663872142,1.000252529336232,0.9997559164634675,1.000222298189604
12|10|OffsetF45|1.008658063000248,0.9776599495876245,1.023129671862989,1.023129671862989,0.992616790466759,1.012725007352464,1.008440267632869,0.9656788526781157,1.023129671862989,0.9656788526781157,0.993046158078209,1.023129671862989,1.008440267632869,0.9870595601554923,1.000528529049814,0.9990829307246578,0.9889119702592479,0.9990829307246578,0.9974663956601831,1.002599050900744,0.9977443450824096
13|11|OffsetA45|1.002206115577125,0.9941891788109908,1.005865936989239,1.005865936989239,0.9980798403250863,1.00332901908819,1.002295542145029,0.991047796534645,1.005865936989239,0.991047796534645,0.9982463486991924,1.005865936989239,1.002295542145029,0.9966716693398372,1.000187691235058,0.999736425222971,0.9971410596682535,0.999736425222971,0.9993222789984999,1.000631010596648,0.9994074722711871
14|12|OffsetFTD|1.003512512846241,0.99039808277526,1.009257961053271,1.009257961053271,0.9968278730740847,1.00555634199889,1.004083142202896,0.9851348086546712,1.009257961053271,0.9851348086546712,0.9972617311560219,1.009257961053271,1.004083142202896,0.9946089754287837,1.00045550763849,0.9994919706369462,0.9953486079772174,0.9994919706369462,0.9988259462462308,1.000912349207639,0.9990040567201671
15|13|OffsetATD|1.000994259708301,0.997610378315901,1.002697115514953,1.002697115514953,0.9992098461689657,1.001332731236358,1.000754173015497,0.9963659284845211,1.002697115514953,0.9963659284845211,0.9991744050892001,1.002697115514953,1.000754173015497,0.998560050750009,0.9999820095922312,0.9999390723383237,0.9987766492229945,0.9999390723383237,0.9997567657959842,1.000344835739778,0.9997673546917121
16|14|Holden2.5|0.9851398444302544,1.03851730482983,0.9603423147990229,0.9603423147990229,1.01272926495976,0.9780321928283118,0.9853004925263786,1.059211141857921,0.9603423147990229,1.059211141857921,1.011908308444843,0.9603423147990229,0.9853004925263786,1.02225581719215,0.9990148787430694,1.001617983061108,1.019080299781353,1.001617983061108,1.004395803535771,0.9955851057884212,1.00389760808454
17|15|Holden5|0.9878431391820179,1.031551057107353,0.9675661425075506,0.9675661425075506,1.010426954710481,0.9819987579953959,0.9879248365056952,1.048510743968409,0.9675661425075506,1.048510743968409,1.009735750001818,0.9675661425075506,0.9879248365056952,1.018217790237052,0.9991758956080883,1.001333907660552,1.015620845061528,1.001333907660552,1.003607169423034,0.9963989640399372,1.00319467274734
18|16|Holden10|0.9932617070093639,1.017586520595366,0.9820455392912335,0.9820455392912335,1.005811746672261,0.989950322632615,0.9931865824112508,1.027060488674374,0.9820455392912335,1.027060488674374,1.00538118447501,0.9820455392912335,0.9931865824112508,1.010123535542812,0.9994991645417561,1.000764203458953,1.008686272585962,1.000764203458953,1.00202608926228,0.9980299792680707,1.001785511470184
19|17|Holden15|0.9926397417047704,1.019225456580995,0.9803919241329497,0.9803919241329497,1.006353330211851,0.9890112644708193,0.9925385514917994,1.029585592193007,0.9803919241329497,1.029585592193007,1.005875484839764,0.9803919241329497,0.9925385514917994,1.011062071842403,0.999445940851825,1.000838689272586,1.009492491298886,1.000838689272586,1.002217349507402,0.9978522644325114,1.001952668203244
20|18|HexGrid-90degTilt2.5degRes|0.9999970622579328,1.000008323260156,0.9999923350855698,0.9999923350855698,1.00000274911158,0.9999951427684768,0.9999962252800662,1.000012944473179,0.9999923350855698,1.000012944473179,1.00000223904552,0.9999923350855698,0.9999962252800662,1.000004584876623,0.9999994871958233,1.000000501612938,1.000003972954942,1.000000501612938,1.000001063487599,0.9999993178071931,1.00000087766031
21|19|HexGrid-90degTilt5degRes|0.9999336028516159,1.000172840622855,0.9998229867315456,0.9998229867315456,1.00005711935886,0.999901307204188,0.9999334128615736,1.000265857808287,0.9998229867315456,1.000265857808287,1.000053087618617,0.9998229867315456,0.9999334128615736,1.000099635334931,0.9999952661102169,1.000007419133802,1.000085463342907,1.000007419133802,1.000019844190067,0.9999804726983624,1.000017526882193
22|20|HexGrid-90degTilt10degRes|0.9997656193592607,1.000608014961426,0.9993746258783137,0.9993746258783137,1.000200937088134,0.9996531494141123,0.9997675354745221,1.000934785028049,0.9993746258783137,1.000934785028049,1.000187739851162,0.9993746258783137,0.9997675354745221,1.000351160251286,0.9999842362813279,1.000025648793628,1.000301085863568,1.000025648793629,1.000069470867255,0.9999305018694666,1.000061550881873
23|21|HexGrid-90degTilt15degRes|0.9994502629097619,1.001423873552513,0.9985326767622211,0.9985326767622211,1.000470568697002,0.9991880838717946,0.9994574690692285,1.002188648437274,0.9985326767622211,1.002188648437274,1.000440688434662,0.9985326767622211,0.9994574690692285,1.000823058753251,0.9999640188831151,1.000059598089575,1.000705562067835,1.000059598089575,1.000162340741432,0.9998364079455895,1.000144033966807
"@

$rowLines = $dataRows -split "`n"
foreach ($line in $rowLines) {
    $line = $line.Trim()
    if ($line.Length -eq 0) { continue }
    $fields = $line -split "\|"
    $rowNum = [int]$fields[0]
    $aVal = [int]$fields[1]
    $bLabel = $fields[2]
    $values = $fields[3] -split ","

    $ws.Cells.Item($rowNum, 1).Value = $aVal
    $ws.Cells.Item($rowNum, 2).Value = $bLabel

    $colIdx = 3
    foreach ($v in $values) {
        $ws.Cells.Item($rowNum, $colIdx).Value = [double]$v
        $colIdx++
    }
}
